$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45177 -> 45178) for every data row (rows 2 through 321).
$ws.Range("C2:C321").Value = 45178
